$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.668637275695801
$ws.Range("B1").Value = 3.252068996429443
$ws.Range("C1").Value = 6.070633888244629
$ws.Range("D1").Value = 1.865052700042725
$ws.Range("E1").Value = 0.9073725342750549
